# add Label (Drug) api
# Rename the Thai column headers to English API-style field names,
# and normalize the Amoxycillin label spacing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: translate Thai labels to API field names
$ws.Range("A1").Value = "drug_name"
$ws.Range("B1").Value = "term_of_use"
$ws.Range("C1").Value = "warning_label"
$ws.Range("D1").Value = "text_label"
$ws.Range("E1").Value = "remark"

# Normalize Amoxycillin labels (add a space before the parenthesis)
$ws.Range("A14").Value = "Amoxycillin (ชนิดเม็ด)"
$ws.Range("A17").Value = "Amoxycillin (ชนิดน้ำ)"

# Update the active selection to reflect where the user last clicked
$ws.Range("B7:B9").Select()
